$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4187560232757783
$ws.Range("C2").Value = 0.5852816327940494
$ws.Range("D2").Value = 0.6767035828922799
$ws.Range("E2").Value = 0.740830947129459
$ws.Range("B3").Value = 0.4435718024875882
$ws.Range("C3").Value = 0.6139259745853569
$ws.Range("D3").Value = 0.6980272502158601
$ws.Range("E3").Value = 0.7557067841118835
$ws.Range("B4").Value = 0.3750959194988736
$ws.Range("C4").Value = 0.5313749586002714
$ws.Range("D4").Value = 0.6244327983843749
$ws.Range("E4").Value = 0.695716271214237
$ws.Range("B5").Value = 0.479423323321474
$ws.Range("C5").Value = 0.6562369809696084
$ws.Range("D5").Value = 0.7258354903946782
$ws.Range("E5").Value = 0.7775079597791049
$ws.Range("B6").Value = 0.4703025743668272
$ws.Range("C6").Value = 0.6451873364608633
$ws.Range("D6").Value = 0.7154820880022621
$ws.Range("E6").Value = 0.7695806858902364
$ws.Range("B7").Value = 0.4933860710767114
$ws.Range("C7").Value = 0.6714093701379333
$ws.Range("D7").Value = 0.7460343174914137
$ws.Range("E7").Value = 0.7866323612145677
$ws.Range("B8").Value = 0.4291178983441817
$ws.Range("C8").Value = 0.5986335783223746
$ws.Range("D8").Value = 0.686666898898007
$ws.Range("E8").Value = 0.7321063414733795
$ws.Range("B9").Value = 0.4972499334998906
$ws.Range("C9").Value = 0.6780818500958263
$ws.Range("D9").Value = 0.7531255498779952
$ws.Range("E9").Value = 0.793596799859918
$ws.Range("B10").Value = 0.5111134471564462
$ws.Range("C10").Value = 0.6863841418024282
$ws.Range("D10").Value = 0.7535477510794063
$ws.Range("E10").Value = 0.7849380575792405
$ws.Range("B11").Value = 0.5063591113775412
$ws.Range("C11").Value = 0.6794344659981938
$ws.Range("D11").Value = 0.7456990317317909
$ws.Range("E11").Value = 0.7772108670126719
$ws.Range("B12").Value = 0.4729543549636881
$ws.Range("C12").Value = 0.6156714985252651
$ws.Range("D12").Value = 0.6540875005989942
$ws.Range("E12").Value = 0.6700888040979659
$ws.Range("B13").Value = 0.508822791869401
$ws.Range("C13").Value = 0.6825026344982311
$ws.Range("D13").Value = 0.7481255078652088
$ws.Range("E13").Value = 0.7789775295704988
